$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.353672031788087, -1.292084699452477, 0.5887890085463383, 0.4292840851827592, 0.7152945399284363, 0.3071393668651581, 0.6648047566413879, 0.4754526019096375)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
